# "last uapdate at all" -- append two new NBR_CASS reference rows (ELINE -
# 09178 / IDB-32-3SM-1-A) below the existing table on Feuil1, keeping the
# same look (number format + border) as the row directly above them, then
# leave the selection where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 (3M BPEO T1 (CDP)) already carries the "bordered number" look we
# want for the two new rows -- clone its formatting down first ...
$ws.Range("A31:C31").Copy($ws.Range("A32:C32"))
$ws.Range("A31:C31").Copy($ws.Range("A33:C33"))

# ... then overwrite the values/labels for the two new entries.
$ws.Range("A32").Value = "ELINE - 09178"
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = 12

$ws.Range("A33").Value = "IDB-32-3SM-1-A"
$ws.Range("B33").Value = 6
$ws.Range("C33").Value = 12

# Match the author's final cursor position/selection.
$ws.Range("B35").Select() | Out-Null
